# Applies the "Add files via upload" revision:
#  - Team ID changes (new PNT id)
#  - Team member names on rows 6-9 changed to the new roster
#  - Active sheet switches from "Testscearnios" to "Shopenzer Testcases"
#  - Selection on the main sheet moves to L8

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Shopenzer Testcases")

# Update Team ID (F2) on the main sheet.
$ws1.Range("F2").Value = "PNT2022TMID53380"

# Update the team-member roster in column N (rows 6-9).
$ws1.Range("N6").Value = "Ritunjay M"
$ws1.Range("N7").Value = "Praveen Raagul R"
$ws1.Range("N8").Value = "Pradeep V"
$ws1.Range("N9").Value = "Munish Kumar S"

# Make "Shopenzer Testcases" the active (selected) tab, then move the
# selection/active cell on it to L8.
$ws1.Activate()
$ws1.Range("L8").Select()
